$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: title - new text, merged A1:I1, taller row
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Keda Municipality"
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# Row 2: text unchanged, but row height reverts to the default (no explicit
# custom height any more)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# Row 3: header row (years) - style of A3 changes (font Arial10 -> Sylfaen11)
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# Row 4: label changes from "Number of disability persons" to
# "family with disabilities Persons " and the placeholder ellipsis values are
# replaced with real numbers
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 865
$ws.Range("C4").Value = 823
$ws.Range("D4").Value = 812
$ws.Range("E4").Value = 822
$ws.Range("F4").Value = 855
$ws.Range("G4").Value = 856
$ws.Range("H4").Value = 864
$ws.Range("I4").Value = 881
$ws.Range("B4:I4").NumberFormat = "# ##0"
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# Row 5: previously held the merged "Source:" note; now becomes a data row
# "disabilities Persons " with its own figures. Unmerge first.
# ---------------------------------------------------------------------------
$ws.Range("A5:H5").UnMerge()
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 1038
$ws.Range("C5").Value = 1000
$ws.Range("D5").Value = 979
$ws.Range("E5").Value = 991
$ws.Range("F5").Value = 1021
$ws.Range("G5").Value = 1017
$ws.Range("H5").Value = 1027
$ws.Range("I5").Value = 1050
$ws.Range("B5:I5").NumberFormat = "# ##0"
$ws.Rows.Item(5).RowHeight = 21

# A5 should now carry the same formatting the old A4 label had (border on top)
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# Row 6: now holds the "Source:" text (merged A6:H6), the old "Note:" text is
# removed entirely
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$ws.Range("A6:H6").Merge()
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# Column width: column A narrower, other columns revert to workbook default
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.81640625

$ws.Range("A1").Select()
